# Generate Report for Archive
#
# The three rows describing 3556bb03-...md, 93843a23-...md and
# b4e0931c-...md are rotated: the "93843a23" file and the "b4e0931c" file
# move up (status becomes "In Translation"), and "3556bb03" moves down
# into the vacated slot (keeping "Ready for handoff"). This is applied as
# a plain text edit of the existing cells (hyperlinks are left attached to
# their original rows, matching how Excel behaves when only the cell text
# is overwritten).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet (File Name / zh-cn / de-de) ----
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A4").Value = "93843a23-7e7a-48d9-8fe0-7a25a915ef9f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

$ws.Range("A5").Value = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "In Translation"

$ws.Range("A6").Value = "3556bb03-c6cc-4123-904b-60db6fd4231a.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "Ready for handoff"

# ---- zh-cn sheet ----
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A4").Value = "93843a23-7e7a-48d9-8fe0-7a25a915ef9f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "93843a23-7e7a-48d9-8fe0-7a25a915ef9f.1c1652b8b7cd6e9c6c53d488f4e5a986866af732.zh-cn.xlf"

$ws.Range("A5").Value = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.zh-cn.xlf"

$ws.Range("A6").Value = "3556bb03-c6cc-4123-904b-60db6fd4231a.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "3556bb03-c6cc-4123-904b-60db6fd4231a.f57cf5fe082fb2dd10368d47c58deba6174c54a7.zh-cn.xlf"

# ---- de-de sheet ----
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A4").Value = "93843a23-7e7a-48d9-8fe0-7a25a915ef9f.md"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "93843a23-7e7a-48d9-8fe0-7a25a915ef9f.1c1652b8b7cd6e9c6c53d488f4e5a986866af732.de-de.xlf"

$ws.Range("A5").Value = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.md"
$ws.Range("B5").Value = "In Translation"
$ws.Range("C5").Value = "b4e0931c-0211-423b-a6a7-1fe3fb71bae9.8261d65499c1c71b25bcb0f4c39c0a00967bb7b1.de-de.xlf"

$ws.Range("A6").Value = "3556bb03-c6cc-4123-904b-60db6fd4231a.md"
$ws.Range("B6").Value = "Ready for handoff"
$ws.Range("C6").Value = "3556bb03-c6cc-4123-904b-60db6fd4231a.f57cf5fe082fb2dd10368d47c58deba6174c54a7.de-de.xlf"
